$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix path typo in absPath (GIt -> Git) - this is workbook-level metadata,
# typically not directly settable via COM; Excel will re-derive this on save
# based on the actual file location, so we leave it to the runtime.

# Highlight certain Part cells in column A (and some F cells) yellow
$yellowRows = @(3,4,5,6,7,10,11,15,16,17,18,19,20,21,22,23,24,25,26,27,31,33)
foreach ($r in $yellowRows) {
    $ws.Cells.Item($r, 1).Interior.Color = 65535
}
$ws.Cells.Item(15, 6).Interior.Color = 65535
$ws.Cells.Item(21, 6).Interior.Color = 65535
$ws.Cells.Item(22, 6).Interior.Color = 65535
$ws.Cells.Item(25, 6).Interior.Color = 65535

# Add new rows 35-37 (Parts section)
$ws.Range("A35").Value = "Parts"
$ws.Cells.Item(35,1).Interior.Color = 65535

$ws.Range("A36").Value = "Power Switch"
$ws.Range("C36").Value = "108-0041-EVX"
$ws.Range("D36").Value = "Mountain Switch"
$ws.Range("E36").Value = "Toggle Switches SPST OFF-ON"

$ws.Range("A37").Value = "Reset Switch"
$ws.Range("C37").Value = "103-1013-EVX"
$ws.Range("D37").Value = "Mountain Switch"
$ws.Range("E37").Value = "Pushbutton Switches METAL BODY BLK"

# Update view: remove topLeftCell, update selection
$ws.Range("D38").Select()

# Page setup - set orientation to portrait
$ws.PageSetup.Orientation = 1
